$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("S4").Copy()
$ws.Range("Y1").PasteSpecial(-4104)
Write-Host "done"
